# Upload mission screenshots and notes from 2021-03-25 testing
#
# Fills in the "Log File Name" column (D) for the 2021-03-25 waypoint /
# simple mission rows (45-52) that were previously left blank, and leaves
# the active selection where the author ended up (cell D52) with the
# window scrolled so row 37 is the top visible row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logNames = @{
    45 = "LOG_PEARL_WAYPOINT_25_3_2021_____13_39_28"
    46 = "LOG_PEARL_WAYPOINT_25_3_2021_____13_52_14"
    47 = "LOG_PEARL_WAYPOINT_25_3_2021_____14_05_23"
    48 = "LOG_PEARL_WAYPOINT_25_3_2021_____14_07_52"
    49 = "LOG_PEARL_WAYPOINT_25_3_2021_____14_19_12"
    50 = "LOG_PEARL_WAYPOINT_25_3_2021_____14_32_25"
    51 = "LOG_PEARL_WAYPOINT_25_3_2021_____14_35_08"
    52 = "LOG_PEARL_SIMPLE_25_3_2021_____14_50_18"
}

foreach ($row in $logNames.Keys | Sort-Object) {
    $ws.Cells.Item($row, 4).Value = $logNames[$row]
}

# Scroll the view and move the selection the way the author left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 2

$ws.Range("D52").Select()
